# Weekly data refresh: insert two new price rows (new report date) at the
# top of the Hortaliza/Lechuga detail block (rows 452-453), pushing all the
# existing rows down by two positions. The used range grows from 561 to 563
# rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 452:453 - this shifts the existing rows 452-561
# down to 454-563 and carries their formatting with them.
$ws.Rows("452:453").Insert()

# --- New row 452 : Lechuga / Escarola / Primera -------------------------
$ws.Cells.Item(452, 1).Value2  = 1
$ws.Cells.Item(452, 2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(452, 3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(452, 4).Value2  = 44543
$ws.Cells.Item(452, 5).Value2  = 15
$ws.Cells.Item(452, 6).Value2  = 100112033
$ws.Cells.Item(452, 7).Value2  = "Lechuga"
$ws.Cells.Item(452, 8).Value2  = "Escarola"
$ws.Cells.Item(452, 9).Value2  = "Primera"
$ws.Cells.Item(452, 10).Value2 = 120
$ws.Cells.Item(452, 11).Value2 = 2000
$ws.Cells.Item(452, 12).Value2 = 2500
$ws.Cells.Item(452, 13).Value2 = 2250
$ws.Cells.Item(452, 14).Value2 = '$/caja 12 unidades'
$ws.Cells.Item(452, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(452, 16).Value2 = 188
$ws.Cells.Item(452, 17).Value2 = 12
$ws.Cells.Item(452, 18).Value2 = "Hortaliza"

# --- New row 453 : Lechuga / Escarola / Segunda --------------------------
$ws.Cells.Item(453, 1).Value2  = 1
$ws.Cells.Item(453, 2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(453, 3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(453, 4).Value2  = 44543
$ws.Cells.Item(453, 5).Value2  = 15
$ws.Cells.Item(453, 6).Value2  = 100112033
$ws.Cells.Item(453, 7).Value2  = "Lechuga"
$ws.Cells.Item(453, 8).Value2  = "Escarola"
$ws.Cells.Item(453, 9).Value2  = "Segunda"
$ws.Cells.Item(453, 10).Value2 = 160
$ws.Cells.Item(453, 11).Value2 = 2000
$ws.Cells.Item(453, 12).Value2 = 2500
$ws.Cells.Item(453, 13).Value2 = 2250
$ws.Cells.Item(453, 14).Value2 = '$/caja 18 unidades'
$ws.Cells.Item(453, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(453, 16).Value2 = 125
$ws.Cells.Item(453, 17).Value2 = 18
$ws.Cells.Item(453, 18).Value2 = "Hortaliza"
